$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Shared-string driven text updates
#    "In Translation" -> "Handed back: in sync with en-US"
#    This string is used by Overview!E2,F2,E3,F3 (and nowhere else), so we
#    simply rewrite those four cells; the engine pools identical strings so
#    this has the same end effect as editing the shared-string table.
# ---------------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = "Handed back: in sync with en-US"
$ovw.Range("F2").Value = "Handed back: in sync with en-US"
$ovw.Range("E3").Value = "Handed back: in sync with en-US"
$ovw.Range("F3").Value = "Handed back: in sync with en-US"

# Overview columns E/F (Status columns) grow to fit the longer text.
$ovw.Columns.Item(5).ColumnWidth = 29.17
$ovw.Columns.Item(6).ColumnWidth = 29.17

# ---------------------------------------------------------------------------
# 2. zh-cn sheet (handback for the zh-cn locale)
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zhUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d510af8550b400fb9e71aaaf4d20de4599e05de/e2e/15f887d1-4278-4a38-8dad-a75a0ea926aa.md"
$zhUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d510af8550b400fb9e71aaaf4d20de4599e05de/e2e/e351aad4-0f86-46b2-8f5c-8cd90d9be0c5.md"

# Fill in "Latest Target File" (I) and "Latest Handback File" (J) for both rows.
$zh.Range("J2").Value = "15f887d1-4278-4a38-8dad-a75a0ea926aa.323bf7dbaca49c858aa3382b4c4f199099ec4ec2.zh-cn.xlf"
$zh.Range("J3").Value = "e351aad4-0f86-46b2-8f5c-8cd90d9be0c5.511458eff13832f681fa374850aac041bda40f20.zh-cn.xlf"

# Latest Handback DateTime (K) for both rows - both files handed back together.
$zh.Range("K2").Value = "2016-08-18 18:25:31"
$zh.Range("K3").Value = "2016-08-18 18:25:31"

# Rebuild the hyperlinks in document order (A2, I2, A3, I3) so the new "Latest
# Target File" links (I2/I3) sit next to the existing source-file links.
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $zhUrl1, "", "", "15f887d1-4278-4a38-8dad-a75a0ea926aa.md")
$zh.Hyperlinks.Add($zh.Range("I2"), $zhUrl1, "", "", "15f887d1-4278-4a38-8dad-a75a0ea926aa.md")
$zh.Hyperlinks.Add($zh.Range("A3"), $zhUrl2, "", "", "e351aad4-0f86-46b2-8f5c-8cd90d9be0c5.md")
$zh.Hyperlinks.Add($zh.Range("I3"), $zhUrl2, "", "", "e351aad4-0f86-46b2-8f5c-8cd90d9be0c5.md")

# Column widths: Status (C) grows, and the two new hyperlink/file columns (I, J)
# widen to fit the long file names.
$zh.Columns.Item(3).ColumnWidth = 29.17
$zh.Columns.Item(9).ColumnWidth = 39.17
$zh.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# 3. de-de sheet (handback for the de-de locale)
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$deUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d510af8550b400fb9e71aaaf4d20de4599e05de/e2e/15f887d1-4278-4a38-8dad-a75a0ea926aa.md"
$deUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d510af8550b400fb9e71aaaf4d20de4599e05de/e2e/e351aad4-0f86-46b2-8f5c-8cd90d9be0c5.md"

# Row 3's "Latest Handoff File" (G) now points at the de-de xlf that was generated.
$de.Range("G3").Value = "e351aad4-0f86-46b2-8f5c-8cd90d9be0c5.511458eff13832f681fa374850aac041bda40f20.de-de.xlf"

# Fill in "Latest Target File" (I) and "Latest Handback File" (J) for both rows.
$de.Range("J2").Value = "15f887d1-4278-4a38-8dad-a75a0ea926aa.323bf7dbaca49c858aa3382b4c4f199099ec4ec2.de-de.xlf"
$de.Range("J3").Value = "e351aad4-0f86-46b2-8f5c-8cd90d9be0c5.511458eff13832f681fa374850aac041bda40f20.de-de.xlf"

# Latest Handback DateTime (K) for both rows - both files handed back together.
$de.Range("K2").Value = "2016-08-18 18:25:39"
$de.Range("K3").Value = "2016-08-18 18:25:39"

# Rebuild the hyperlinks in document order (A2, I2, A3, I3).
$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $deUrl1, "", "", "15f887d1-4278-4a38-8dad-a75a0ea926aa.md")
$de.Hyperlinks.Add($de.Range("I2"), $deUrl1, "", "", "15f887d1-4278-4a38-8dad-a75a0ea926aa.md")
$de.Hyperlinks.Add($de.Range("A3"), $deUrl2, "", "", "e351aad4-0f86-46b2-8f5c-8cd90d9be0c5.md")
$de.Hyperlinks.Add($de.Range("I3"), $deUrl2, "", "", "e351aad4-0f86-46b2-8f5c-8cd90d9be0c5.md")

# Column widths: Status (C) grows, and the two new hyperlink/file columns (I, J)
# widen to fit the long file names.
$de.Columns.Item(3).ColumnWidth = 29.17
$de.Columns.Item(9).ColumnWidth = 39.17
$de.Columns.Item(10).ColumnWidth = 39.17
